# overlap_percentage.xlsx — update formulas to compute IOP (Image Overlap
# Percentage): the percentage of overlap between two consecutive frames
# used in the VO estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# INPUTS block (J:N) — new camera-height / pitch values
# ---------------------------------------------------------------------
$ws.Range("K3").Value = 0.633
$ws.Range("K4").Value = 30
$ws.Range("K7").Value = 0.07
$ws.Range("K8").Value = 1

# ---------------------------------------------------------------------
# CALCULATIONS block (A:C) — updated formulas
# ---------------------------------------------------------------------
$ws.Range("B9").Formula = "=MAX(0,B7-K8*K7*K11)"

$ws.Range("A10").Value = "theta"

$ws.Range("B11").Formula = "=B5+2*B12*TAN(B10)/(B4-B3)"

# new row: l overlap
$ws.Range("A12").Value = "l overlap"
$ws.Range("B12").Formula = "=B9/SIN(PI()/2 -B10)"
$ws.Range("C12").Value = "m"

# IOP (was "overlap percentage"), now expressed as a true percentage (*100)
$ws.Range("A26").Value = "IOP"
$ws.Range("B26").Formula = "=((B8+B11)*B9/2)/((B5+B6)*B7/2)*100"

# ---------------------------------------------------------------------
# Secondary reference inputs (P:Q) — LOCCAM / NAVCAM comparison values
# ---------------------------------------------------------------------
$ws.Range("P2").Value = "LOCCAM"
$ws.Range("Q2").Value = "NAVCAM"
$ws.Range("P3").Value = 0.24
$ws.Range("Q3").Value = 0.633
$ws.Range("P4").Value = 30.4

# ---------------------------------------------------------------------
# IFD / IOP lookup table used for the chart (H19:J31)
# ---------------------------------------------------------------------
$ws.Range("H19").Value = "CAMERA"
$ws.Range("I19").Value = "loccam"
$ws.Range("J19").Value = "loccam"

$ws.Range("H20").Value = "PITCH"
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 30

$ws.Range("H21").Value = "SPEED"
$ws.Range("I21").Value = "0.07m/s"
$ws.Range("J21").Value = "0.02m/s"

$ws.Range("H22").Value = "IFD"

$ws.Range("H23").Value = 0.01

$ws.Range("H24").Value = 0.03
$ws.Range("J24").Value = 0.92340458663566416

$ws.Range("H25").Value = 0.06
$ws.Range("J25").Value = 0.89903374961988936

$ws.Range("H26").Value = 0.1
$ws.Range("I26").Value = 0.86704400000000004
$ws.Range("J26").Value = 0.86704400000000004

$ws.Range("A27").Value = "IFD"
$ws.Range("B27").Formula = "=K7*K8"
$ws.Range("H27").Value = 0.2
$ws.Range("J27").Value = 0.789594989642042

$ws.Range("H28").Value = 0.3
$ws.Range("J28").Value = 0.71575232521662158

$ws.Range("H29").Value = 0.5
$ws.Range("J29").Value = 0.57888669240552937

$ws.Range("H30").Value = 0.75
$ws.Range("J30").Value = 0.42809158146619208

$ws.Range("H31").Value = 1
$ws.Range("J31").Value = 0.2998375039429973

# ---------------------------------------------------------------------
# Column width for the new lookup-table column (J)
# ---------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 9.85546875

# ---------------------------------------------------------------------
# View state — scroll so row 2 is pinned at the top, select C13
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$win.Width = 874.5
$win.Height = 394.5
$ws.Range("C13").Select()
